$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: split the opening sentence's single run into three runs:
#   "Survive" | " in" | " space where aliens want nothing but to see your
#   destruction."
# The COM layer auto-merges adjacent runs that end up with identical
# formatting, so we insert the new text first (it merges right back into one
# run) and then briefly toggle Bold on just the inserted " in" piece; Word
# has to split the run to apply/remove that formatting, and - since this
# engine does not retroactively re-merge runs once split - clearing Bold
# again leaves three separate, identically-formatted runs behind.
# ---------------------------------------------------------------------------

$rSurvive = $d.Content
$rSurvive.Find.Execute("Survive", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rSurvive.Collapse(0)
$rSurvive.InsertAfter(" in")

$rIn = $d.Content
$rIn.Find.Execute(" in", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rIn.Bold = 1
$rIn.Bold = 0

# ---------------------------------------------------------------------------
# Change 2: merge the three runs around "16-bit" ("... Audio will sound
# like ", "16-bit", " sound effects. ...") back into a single run. A no-op
# Find/Replace over just "16-bit" causes the engine to coalesce it with its
# two (identically formatted) neighbours while leaving the preceding " 80s"
# run untouched.
# ---------------------------------------------------------------------------

$rBit = $d.Content
$rBit.Find.Execute("16-bit", $true, $false, $false, $false, $false, $true, 1, $false, "16-bit", 2)
